$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.581.32"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.730.92"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.45"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2677"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06191"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.731.23"
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07189"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.61"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6118"
$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.539"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.31"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9992"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.583.13"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006977"
$ws.Range("E19").Value = "  +1.14%  "

$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.952.27"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.535"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.823"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.267"
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.36"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.400"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.54"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.991"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08034"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.702"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04530"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.611"
$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.008"
$ws.Range("E35").Value = "  +1.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6283"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.073"
$ws.Range("E37").Value = "  +5.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9074"
$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.412"
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01508"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.36"
$ws.Range("E42").Value = "  -11.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.514"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3904"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.061"
$ws.Range("E45").Value = "  +7.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1184"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05384"

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.885"
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.76"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.251"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3421"
$ws.Range("E51").Value = "  -0.37%  "
